$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column W header: copy formatting from V1 (bold/centered/bordered header
# style) then set the new shared-string text, mirroring how the existing
# header cells (U1, V1, ...) were built.
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)
$ws.Range("W1").Value = "param_E_pv3_solar"

# Existing column U ("param_P_to_charging_station1") data values change
# from 0 to 0.12 for every data row.
$ws.Range("U2:U6").Value = 0.12

# New column W data values are 0 for every data row.
$ws.Range("W2:W6").Value = 0
